# StructureDefinition-crossover-indicator.xlsx
# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Changes:
#   Metadata sheet:
#     - Version bumped 5.0.0 -> 6.0.0
#     - Date bumped to the new publication timestamp
#     - Publisher's "Contact" rows (a ContactDetail with no usable display
#       text, duplicated on two rows) are replaced by a single
#       "Jurisdiction" / "United States of America" row, and the
#       "Publisher" value is now populated ("Alvearie Team").
#   Elements sheet:
#     - The root Extension row's Short/Definition no longer show the
#       generic "Extension" / "An Extension" placeholder text; they now
#       show the real title/description of this extension.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Elements"

# --- Metadata sheet -------------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date: bump to the new publish timestamp
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$ws1.Range("B9").Value = "Alvearie Team"

# First "Contact" row becomes "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Second duplicate "Contact" / "No display for ContactDetail" row is removed
# entirely, shifting Description..Context up by one row.
$ws1.Rows("11:11").Delete()

# --- Elements sheet ---------------------------------------------------

# Root Extension row: Short/Definition get the extension's real title and
# description instead of the generic placeholder text.
$ws2.Range("K2").Value = "Crossover Indicator"
$ws2.Range("L2").Value = "Indicates whether the claim is a crossover claim where a portion is paid by Medicare"
